$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Update the username embedded in a few text values from user4 to user2
$ws.Range("B1").Value2 = "template /pub/home/user2/jmt_workspace/blocks/bl_1s13/bl_1s13.tsdl"
$ws.Range("K3").Value2 = "/pub/home/user2/jmt_workspace/workshop_config.sdl"
$ws.Range("J3").Value2 = "/pub/home/user2/jmt_workspace"

# Update the active selection on the sheet
$ws.Range("J5").Select()

# Maximize the workbook window (reflects the new saved window size/position)
$excel.ActiveWindow.WindowState = -4137
